$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 5 year value from 2030 to 2025
$ws.Range("F5").Value = 2025

# Add 5 more rows (rows 6-10) that replicate row 5's content/format,
# with Year values 2030, 2035, 2040, 2045, 2050
$years = @(2030, 2035, 2040, 2045, 2050)
$destRow = 6
foreach ($year in $years) {
    $srcRange = $ws.Range("B5:H5")
    $dstRange = $ws.Range("B" + $destRow + ":H" + $destRow)
    $srcRange.Copy($dstRange)
    $ws.Range("F" + $destRow).Value = $year
    $destRow++
}

# Match the new active selection recorded in the workbook after the edit
$ws.Range("I15").Select()
